$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CenDLL")
$ws.Range("B3").Value = "test"
